$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the task name for the "training/driving" row to reflect the new clutter-change task
$ws.Range("B6").Value = "full_task_training"

# Update the active selection to match the latest edit location
$ws.Range("B7").Select()
